$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsProductos   = $wb.Worksheets.Item("Productos")
$wsUsuarios    = $wb.Worksheets.Item("Usuarios")
$wsCategorias  = $wb.Worksheets.Item("Categorias")
$wsProveedores = $wb.Worksheets.Item("Proveedores")
$wsHistorial   = $wb.Worksheets.Item("HistorialVentas")

# 1) Productos: duplicate the last product row (row 13) into a new row 14
$wsProductos.Range("A13:J13").Copy($wsProductos.Range("A14:J14"))

# 2) Categorias currently only has its header row; the "Categorias" rows were
#    mistakenly saved into the Usuarios sheet. Move that data (A2:B6) over to
#    the Categorias sheet, below its header.
$wsUsuarios.Range("A2:B6").Copy($wsCategorias.Range("A2:B6"))
$wsCategorias.Columns("A").ColumnWidth = 14.166666666666666
$wsCategorias.Columns("B").ColumnWidth = 20.830729166666668

# 3) Usuarios should instead hold the real user records. Clear out the
#    mistaken categoria rows and write the correct Usuario data.
$wsUsuarios.Range("A2:D6").ClearContents()

$wsUsuarios.Range("A2").Value = "USR-001"
$wsUsuarios.Range("B2").Value = "admin"
$wsUsuarios.Range("C2").Value = "admin123"
$wsUsuarios.Range("D2").Value = "Administrador"

$wsUsuarios.Range("A3").Value = "USR-002"
$wsUsuarios.Range("B3").Value = "empleado1"
$wsUsuarios.Range("C3").Value = "emp456"
$wsUsuarios.Range("D3").Value = "Empleado"

$wsUsuarios.Range("A4").Value = "USR-003"
$wsUsuarios.Range("B4").Value = "jefe"
$wsUsuarios.Range("C4").Value = "jefe1"
$wsUsuarios.Range("D4").Value = "Administrador"

# 4) Restore/update the selection on each sheet
$wsUsuarios.Range("C26").Select()
$wsCategorias.Range("C8").Select()
$wsHistorial.Range("E27").Select()

# 5) Make Proveedores the active sheet/tab (was HistorialVentas before)
$wsProveedores.Activate()
$wsProveedores.Range("A2").Select()
